# Update "想去人数" (want-to-go count) values in column F across sheets
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 859
$ws1.Range("F3").Value = 13786
$ws1.Range("F4").Value = 13574
$ws1.Range("F8").Value = 595
$ws1.Range("F13").Value = 2144
$ws1.Range("F17").Value = 124
$ws1.Range("F21").Value = 398
$ws1.Range("F24").Value = 833
$ws1.Range("F25").Value = 83

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 1500

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 859
$ws4.Range("F4").Value = 13786
$ws4.Range("F5").Value = 13574
$ws4.Range("F9").Value = 595
$ws4.Range("F16").Value = 2144
$ws4.Range("F20").Value = 124
$ws4.Range("F28").Value = 398
$ws4.Range("F31").Value = 833
$ws4.Range("F33").Value = 1500
$ws4.Range("F37").Value = 83
